$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = -8.115
$ws.Range("C7").Value = -13.365
$ws.Range("B9").Value = 5.369999999999999
$ws.Range("C12").Value = -11.141
$ws.Range("B13").Value = 5.566999999999999
$ws.Range("C14").Value = -13.321
$ws.Range("D15").Value = -8.084
$ws.Range("B16").Value = 5.257000000000001
$ws.Range("B18").Value = 5.029000000000001
$ws.Range("C19").Value = -12.028
$ws.Range("B20").Value = 6.775999999999999
$ws.Range("B26").Value = 5.767
$ws.Range("C26").Value = -13.358
$ws.Range("B27").Value = 5.577
$ws.Range("C27").Value = -13.625
$ws.Range("D28").Value = -8.096
$ws.Range("B29").Value = 5.355
$ws.Range("C29").Value = -11.281
$ws.Range("D33").Value = -7.144999999999999
$ws.Range("B35").Value = 9.239000000000001
$ws.Range("D35").Value = -7.826000000000001
$ws.Range("B36").Value = 7.971999999999999
$ws.Range("C37").Value = -13.038
$ws.Range("C38").Value = -14.03
$ws.Range("D38").Value = -8.068999999999999
$ws.Range("D43").Value = -7.776999999999999
$ws.Range("D44").Value = -7.476999999999999
$ws.Range("B45").Value = 6.224000000000001
$ws.Range("D45").Value = -7.666000000000001
$ws.Range("C47").Value = -12.094
$ws.Range("D47").Value = -7.341000000000001
$ws.Range("C51").Value = -12.369
$ws.Range("D51").Value = -7.542999999999999
$ws.Range("C52").Value = -11.423
$ws.Range("D54").Value = -8.339
$ws.Range("B55").Value = 5.757
$ws.Range("C55").Value = -13.438
$ws.Range("B57").Value = 5.424
$ws.Range("D57").Value = -8.099
$ws.Range("D62").Value = -7.867
$ws.Range("D63").Value = -7.613999999999999
$ws.Range("D67").Value = -6.9
$ws.Range("B69").Value = 5.3
$ws.Range("C69").Value = -10.921
$ws.Range("C70").Value = -12.195
$ws.Range("D70").Value = -7.668000000000001
$ws.Range("B76").Value = 5.872
$ws.Range("C76").Value = -12.395
$ws.Range("B78").Value = 7.523999999999999
$ws.Range("C81").Value = -12.712
$ws.Range("D81").Value = -7.461
$ws.Range("B82").Value = 5.467000000000001
$ws.Range("B83").Value = 5.351
$ws.Range("C83").Value = -13.97
$ws.Range("D88").Value = -7.651999999999998
$ws.Range("B93").Value = 5.516
$ws.Range("C94").Value = -11.012
$ws.Range("D96").Value = -7.545
$ws.Range("B97").Value = 6.449000000000001
$ws.Range("D99").Value = -7.74
$ws.Range("C100").Value = -13.409
$ws.Range("C102").Value = -13.311
